# Apply scheduled-runner market data refresh to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 914.2857
$ws.Range("I18").Value = 914.2857
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 914.2857
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -630.2857
$ws.Range("H28").Value = 7457008
$ws.Range("I28").Value = 10439324
$ws.Range("J28").Value = 1218
$ws.Range("K28").Value = 10439324
$ws.Range("L28").Value = 1218
$ws.Range("M28").Value = -10438839
$ws.Range("N28").Value = -2188
$ws.Range("H33").Value = 5373.095
$ws.Range("I33").Value = 148.75
$ws.Range("J33").Value = 12338.889
$ws.Range("K33").Value = 148.75
$ws.Range("L33").Value = 12338.889
$ws.Range("M33").Value = 80.25
$ws.Range("N33").Value = -12796.889
$ws.Range("H62").Value = 333334660
$ws.Range("I62").Value = 500001000
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 500001000
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -500000376
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 333334660
$ws.Range("I65").Value = 500001000
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 2500005000
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -2500001880
$ws.Range("N65").Value = -16240
$ws.Range("H113").Value = 2364.3125
$ws.Range("I113").Value = 2338.2144
$ws.Range("J113").Value = 2547
$ws.Range("K113").Value = 2338.2144
$ws.Range("L113").Value = 2547
$ws.Range("M113").Value = 915.7856000000002
$ws.Range("N113").Value = -9055
$ws.Range("H141").Value = 4737.9443
$ws.Range("I141").Value = 4840.1763
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 14520.5289
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -9340.528900000001
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4596.6665
$ws.Range("I63").Value = 2895
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 2895
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -2209
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 4596.6665
$ws.Range("I66").Value = 2895
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 14475
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -11043
$ws.Range("N66").Value = -46864
$ws.Range("H97").Value = 685.2105
$ws.Range("I97").Value = 544.25
$ws.Range("J97").Value = 1437
$ws.Range("K97").Value = 544.25
$ws.Range("L97").Value = 1437
$ws.Range("M97").Value = -48.25
$ws.Range("N97").Value = -2429
$ws.Range("H102").Value = 3326
$ws.Range("I102").Value = 2722.8572
$ws.Range("J102").Value = 4733.3335
$ws.Range("K102").Value = 2722.8572
$ws.Range("L102").Value = 4733.3335
$ws.Range("M102").Value = -1100.8572
$ws.Range("N102").Value = -7977.3335
$ws.Range("H122").Value = 19280126
$ws.Range("I122").Value = 88963.57000000001
$ws.Range("J122").Value = 41669816
$ws.Range("K122").Value = 266890.71
$ws.Range("L122").Value = 125009448
$ws.Range("M122").Value = -264440.71
$ws.Range("N122").Value = -125014348
$ws.Range("H132").Value = 10022009
$ws.Range("I132").Value = 16130827
$ws.Range("J132").Value = 54988.367
$ws.Range("K132").Value = 48392481
$ws.Range("L132").Value = 164965.101
$ws.Range("M132").Value = -48389951
$ws.Range("N132").Value = -170025.101

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 27650.4
$ws.Range("I107").Value = 52762.2
$ws.Range("J107").Value = 2538.6
$ws.Range("K107").Value = 52762.2
$ws.Range("L107").Value = 2538.6
$ws.Range("M107").Value = -50842.2
$ws.Range("N107").Value = -6378.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1253507
$ws.Range("I31").Value = 2779720
$ws.Range("J31").Value = 4787.273
$ws.Range("K31").Value = 2779720
$ws.Range("L31").Value = 4787.273
$ws.Range("M31").Value = -2779425
$ws.Range("N31").Value = -5377.273
$ws.Range("H34").Value = 1253507
$ws.Range("I34").Value = 2779720
$ws.Range("J34").Value = 4787.273
$ws.Range("K34").Value = 2779720
$ws.Range("L34").Value = 4787.273
$ws.Range("M34").Value = -2779518
$ws.Range("N34").Value = -5191.273
$ws.Range("H39").Value = 1646.8182
$ws.Range("I39").Value = 1171.5
$ws.Range("J39").Value = 6400
$ws.Range("K39").Value = 1171.5
$ws.Range("L39").Value = 6400
$ws.Range("M39").Value = -780.5
$ws.Range("N39").Value = -7182
$ws.Range("H49").Value = 1646.8182
$ws.Range("I49").Value = 1171.5
$ws.Range("J49").Value = 6400
$ws.Range("K49").Value = 1171.5
$ws.Range("L49").Value = 6400
$ws.Range("M49").Value = -989.5
$ws.Range("N49").Value = -6764
$ws.Range("H99").Value = 406518.1
$ws.Range("I99").Value = 516898.2
$ws.Range("J99").Value = 1791.1111
$ws.Range("K99").Value = 516898.2
$ws.Range("L99").Value = 1791.1111
$ws.Range("M99").Value = -515400.2
$ws.Range("N99").Value = -4787.1111
$ws.Range("H105").Value = 1131.7307
$ws.Range("I105").Value = 729.58826
$ws.Range("J105").Value = 1891.3334
$ws.Range("K105").Value = 729.58826
$ws.Range("L105").Value = 1891.3334
$ws.Range("M105").Value = 1017.41174
$ws.Range("N105").Value = -5385.3334
$ws.Range("H126").Value = 406518.1
$ws.Range("I126").Value = 516898.2
$ws.Range("J126").Value = 1791.1111
$ws.Range("K126").Value = 1550694.6
$ws.Range("L126").Value = 5373.3333
$ws.Range("M126").Value = -1548224.6
$ws.Range("N126").Value = -10313.3333
$ws.Range("H134").Value = 18006174
$ws.Range("I134").Value = 20455070
$ws.Range("J134").Value = 47600
$ws.Range("K134").Value = 61365210
$ws.Range("L134").Value = 142800
$ws.Range("M134").Value = -61362675
$ws.Range("N134").Value = -147870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 603.1613
$ws.Range("I107").Value = 361.8095
$ws.Range("J107").Value = 1110
$ws.Range("K107").Value = 1085.4285
$ws.Range("L107").Value = 3330
$ws.Range("M107").Value = 834.5715
$ws.Range("N107").Value = -7170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H107").Value = 731.0769
$ws.Range("I107").Value = 895.9
$ws.Range("J107").Value = 181.66667
$ws.Range("K107").Value = 895.9
$ws.Range("L107").Value = 181.66667
$ws.Range("M107").Value = 1024.1
$ws.Range("N107").Value = -4021.66667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2348.3057
$ws.Range("I100").Value = 1841.7059
$ws.Range("J100").Value = 2801.5789
$ws.Range("K100").Value = 1841.7059
$ws.Range("L100").Value = 2801.5789
$ws.Range("M100").Value = -1300.7059
$ws.Range("N100").Value = -3883.5789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 4260.25
$ws.Range("I39").Value = 2044
$ws.Range("J39").Value = 4999
$ws.Range("K39").Value = 2044
$ws.Range("L39").Value = 4999
$ws.Range("M39").Value = -1631
$ws.Range("N39").Value = -5825
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H126").Value = 1580.8
$ws.Range("I126").Value = 1728.375
$ws.Range("J126").Value = 1412.1428
$ws.Range("K126").Value = 5185.125
$ws.Range("L126").Value = 4236.428400000001
$ws.Range("M126").Value = -2715.125
$ws.Range("N126").Value = -9176.428400000001
$ws.Range("H136").Value = 36509
$ws.Range("I136").Value = 100404.8
$ws.Range("J136").Value = 1011.3333
$ws.Range("K136").Value = 301214.4
$ws.Range("L136").Value = 3033.9999
$ws.Range("M136").Value = -298664.4
$ws.Range("N136").Value = -8133.9999

# A few rows now compute to fully empty trailing cells (no remaining profit figure);
# clear them so the saved cells match the refreshed source data exactly.
$wb.Worksheets.Item("GSM").Range("N18").ClearContents()
$wb.Worksheets.Item("WVR").Range("M122").ClearContents()
$wb.Worksheets.Item("WVR").Range("N122").ClearContents()
